# Addition of VnV feedback partners to Repos.xlsx file
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column K header in row 2 (matches style of existing J2 header cell)
$ws.Range("K2").Value = "VnV Plan"
$ws.Range("K2").Font.Bold = $true

# Set column J width (the newly added <col> entry covers column 10 = J).
# Excel's COM layer round-trips ColumnWidth through pixel metrics, adding a
# small constant offset before it is written back out as the stored
# <col width="..."> value, so compensate to land on the target width (17).
$ws.Columns.Item(10).ColumnWidth = 16.1667

# Fill K3:K14 with formulas that rotate the J column the same way J rotates A
$ws.Range("K3").Formula = "=J4"
$ws.Range("K4").Formula = "=J5"
$ws.Range("K5").Formula = "=J6"
$ws.Range("K6").Formula = "=J7"
$ws.Range("K7").Formula = "=J8"
$ws.Range("K8").Formula = "=J9"
$ws.Range("K9").Formula = "=J10"
$ws.Range("K10").Formula = "=J11"
$ws.Range("K11").Formula = "=J12"
$ws.Range("K12").Formula = "=J13"
$ws.Range("K13").Formula = "=J14"
$ws.Range("K14").Formula = "=J3"

# Update selection to match the post-edit cursor position
$ws.Range("K15").Select()
